# Judge Gift Cards 2008 (G08): re-header the sheet and drop the sample
# card rows, leaving a single blank row below the (new) header, per the
# "Script to download new sets finished" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 becomes the set header/title instead of the generic "Judge" label.
$ws.Range("A1").Value = "Judge Gift Cards 2008 (G08)"

# Remove the four sample card rows (3-6) entirely - delete whole rows so
# the remaining rows below (none, here) shift up cleanly without
# disturbing column A's single-column layout.
$ws.Rows("3:6").Delete()

# Row 2 keeps existing (it is still part of the sheet's dimension) but its
# former "Demonic Tutor" content is wiped out, leaving a blank cell.
$ws.Range("A2").Value = ""
$ws.Range("A2").Style = "Normal"
